$d = $word.ActiveDocument

# --- 1. Insert new "Set for Team 163" heading paragraph at the very start ---
$r = $d.Range(0, 0)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr>' +
  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/>' +
  '</w:rPr></w:pPr>' +
  '<w:r><w:rPr>' +
  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/>' +
  '</w:rPr><w:t>Set for Team 16</w:t></w:r>' +
  '<w:r><w:rPr>' +
  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/>' +
  '</w:rPr><w:t>3</w:t></w:r>' +
  '</w:p>'
$r.InsertXML($newParaXml)

# --- 2. Merge the split "Josephus ... Yodfat ..." runs (and drop proofErr marks) into one run ---
$joseParaIdx = 2
$joseP = $d.Paragraphs($joseParaIdx)
$joseRng = $d.Range($joseP.Range.Start, $joseP.Range.End)
$joseText = "The Josephus' problem is known because of the Flavius Josephus' legend, a Jewish historian living in the 1st century. According to Josephus' account of the siege of Yodfat, he and his 40 comrade soldiers were trapped in a cave, the exit of which one was blocked by Romans. They chose suicide over capture and decided that they would form a circle and start killing themselves skipping three in three. Josephus says that, by luck or maybe by the hand of God, he remained the last and gave up to the Romans.”"
$joseXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>' + $joseText + '</w:t></w:r></w:p>'
$joseRng.InsertXML($joseXml)

# --- 3. Merge the split "First, your program..." runs into one run (leave the trailing run intact) ---
$firstParaIdx = 4
$firstP = $d.Paragraphs($firstParaIdx)
$mergeRng = $d.Range($firstP.Range.Start, $firstP.Range.End)
$mergedText = "First, your program should take input the number of test cases from user. If there are NC (1 ≤ NC ≤ 20) test cases. In each input test case, there will be a pair of positive integer numbers n (1 ≤ n ≤ 500) and k (1 ≤ k ≤ 50). The number n represents the quantity of people in the circle, numbered from 1 to n. The number k represents the size of step between two men in the circle.  The image below illustrates "
$tailText = "an example with 5 men and step 2: In this example the remaining element is 3."
$mergeXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">' + $mergedText + '</w:t></w:r><w:r><w:t>' + $tailText + '</w:t></w:r></w:p>'
$mergeRng.InsertXML($mergeXml)

Write-Output "Paragraph1: $($d.Paragraphs(1).Range.Text)"
Write-Output "Paragraph2: $($d.Paragraphs(2).Range.Text)"
Write-Output "Paragraph3: $($d.Paragraphs(3).Range.Text)"
Write-Output "Paragraph4: $($d.Paragraphs(4).Range.Text)"
Write-Output "TotalParas: $($d.Paragraphs.Count)"
